$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-17 Friday", "2024-05-18 Saturday"),
    @("99×26=", "79×22="),
    @("39×89=", "36×45="),
    @("71×82=", "25×81="),
    @("71×36=", "78×11="),
    @("27×12=", "73×19="),
    @("90×22=", "99×76="),
    @("15×82=", "89×47="),
    @("19×43=", "85×85="),
    @("84×62=", "96×95="),
    @("34×20=", "94×71="),
    @("22×63=", "76×75="),
    @("41×20=", "49×25="),
    @("66×37=", "92×22="),
    @("35×22=", "50×25="),
    @("87×60=", "77×29="),
    @("39×22=", "71×44="),
    @("44×90=", "97×31="),
    @("53×11=", "50×37="),
    @("26×64=", "53×83="),
    @("11×57=", "16×31="),
    @("21×85=", "16×29="),
    @("86×27=", "64×78="),
    @("50×81=", "21×66="),
    @("42×35=", "16×55="),
    @("12×36=", "78×44=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
